$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 179716
$ws.Range("E2").Value = 8781
$ws.Range("F2").Value = 8781
$ws.Range("G2").Value = -2051
$ws.Range("H2").Value = -855
$ws.Range("I2").Value = -947
$ws.Range("J2").Value = 92
$ws.Range("K2").Value = 275519
$ws.Range("L2").Value = 198749
$ws.Range("M2").Value = 76770
$ws.Range("N2").Value = 46594
$ws.Range("O2").Value = 30176
$ws.Range("P2").Value = 5968
$ws.Range("Q2").Value = 6038
$ws.Range("R2").Value = -5447
$ws.Range("S2").Value = 3263
$ws.Range("T2").Value = 3727
$ws.Range("U2").Value = 2311
$ws.Range("V2").Value = 116034
$ws.Range("W2").Value = 4.89
$ws.Range("X2").Value = -0.48
$ws.Range("Y2").Value = -2.02
$ws.Range("Z2").Value = -0.31
$ws.Range("AA2").Value = 258.89
$ws.Range("AB2").Value = 660.76
$ws.Range("AC2").Value = -763
$ws.Range("AD2").Value = -26.73
$ws.Range("AE2").Value = 36136
$ws.Range("AF2").Value = 0.5600000000000001
$ws.Range("AG2").Value = 647
$ws.Range("AH2").Value = 3.17
$ws.Range("AI2").Value = -91.31
$ws.Range("AJ2").Value = 123049283

# Row 3
$ws.Range("D3").Value = 144706
$ws.Range("E3").Value = -273
$ws.Range("F3").Value = -273
$ws.Range("G3").Value = -16127
$ws.Range("H3").Value = -17509
$ws.Range("I3").Value = -10385
$ws.Range("J3").Value = -7124
$ws.Range("K3").Value = 272601
$ws.Range("L3").Value = 202340
$ws.Range("M3").Value = 70261
$ws.Range("N3").Value = 35807
$ws.Range("O3").Value = 34454
$ws.Range("P3").Value = 5968
$ws.Range("Q3").Value = -744
$ws.Range("R3").Value = -3797
$ws.Range("S3").Value = 10586
$ws.Range("T3").Value = 3887
$ws.Range("U3").Value = -4631
$ws.Range("V3").Value = 124201
$ws.Range("W3").Value = -0.19
$ws.Range("X3").Value = -12.1
$ws.Range("Y3").Value = -25.21
$ws.Range("Z3").Value = -6.39
$ws.Range("AA3").Value = 287.98
$ws.Range("AB3").Value = 431.8
$ws.Range("AC3").Value = -7622
$ws.Range("AD3").Value = -2.33
$ws.Range("AE3").Value = 27770
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 733
$ws.Range("AH3").Value = 4.13
$ws.Range("AI3").Value = -9.279999999999999
$ws.Range("AJ3").Value = 123049283

# Row 4
$ws.Range("D4").Value = 139523
$ws.Range("E4").Value = 7982
$ws.Range("F4").Value = 7982
$ws.Range("G4").Value = -2443
$ws.Range("H4").Value = -2155
$ws.Range("I4").Value = -1708
$ws.Range("J4").Value = -448
$ws.Range("K4").Value = 248326
$ws.Range("L4").Value = 180097
$ws.Range("M4").Value = 68228
$ws.Range("N4").Value = 34757
$ws.Range("O4").Value = 33472
$ws.Range("P4").Value = 5968
$ws.Range("Q4").Value = 9676
$ws.Range("R4").Value = 8048
$ws.Range("S4").Value = -22516
$ws.Range("T4").Value = 2909
$ws.Range("U4").Value = 6767
$ws.Range("V4").Value = 101603
$ws.Range("W4").Value = 5.72
$ws.Range("X4").Value = -1.54
$ws.Range("Y4").Value = -4.84
$ws.Range("Z4").Value = -0.83
$ws.Range("AA4").Value = 263.96
$ws.Range("AB4").Value = 394.99
$ws.Range("AC4").Value = -1253
$ws.Range("AD4").Value = -18.73
$ws.Range("AE4").Value = 25509
$ws.Range("AF4").Value = 0.92
$ws.Range("AG4").Value = 474
$ws.Range("AH4").Value = 2.02
$ws.Range("AI4").Value = -41.4
$ws.Range("AJ4").Value = 123049283

# Row 5
$ws.Range("D5").Value = 138413
$ws.Range("E5").Value = 9134
$ws.Range("F5").Value = 9134
$ws.Range("G5").Value = 1309
$ws.Range("H5").Value = -1097
$ws.Range("I5").Value = -2920
$ws.Range("J5").Value = 1824
$ws.Range("K5").Value = 249623
$ws.Range("L5").Value = 183964
$ws.Range("M5").Value = 65659
$ws.Range("N5").Value = 33868
$ws.Range("O5").Value = 31791
$ws.Range("P5").Value = 5968
$ws.Range("Q5").Value = 4293
$ws.Range("R5").Value = -4571
$ws.Range("S5").Value = 7036
$ws.Range("T5").Value = 2773
$ws.Range("U5").Value = 1521
$ws.Range("V5").Value = 108928
$ws.Range("W5").Value = 6.6
$ws.Range("X5").Value = -0.79
$ws.Range("Y5").Value = -8.51
$ws.Range("Z5").Value = -0.44
$ws.Range("AA5").Value = 280.18
$ws.Range("AB5").Value = 352.86
$ws.Range("AC5").Value = -2143
$ws.Range("AD5").Value = -6.18
$ws.Range("AE5").Value = 24847
$ws.Range("AF5").Value = 0.53
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = -5.99
$ws.Range("AJ5").Value = 123402623

# Row 6
$ws.Range("D6").Value = 147611
$ws.Range("E6").Value = 10017
$ws.Range("F6").Value = 10017
$ws.Range("G6").Value = -1052
$ws.Range("H6").Value = -4217
$ws.Range("I6").Value = -5238
$ws.Range("K6").Value = 248149
$ws.Range("L6").Value = 185965
$ws.Range("M6").Value = 62184
$ws.Range("N6").Value = 28048
$ws.Range("P6").Value = 6503
$ws.Range("Q6").Value = 9896
$ws.Range("R6").Value = -7930
$ws.Range("S6").Value = -1101
$ws.Range("T6").Value = 2314
$ws.Range("U6").Value = 7581
$ws.Range("V6").Value = 104784
$ws.Range("W6").Value = 6.79
$ws.Range("X6").Value = -2.86
$ws.Range("Y6").Value = -16.92
$ws.Range("Z6").Value = -1.69
$ws.Range("AA6").Value = 299.05
$ws.Range("AB6").Value = 235.29
$ws.Range("AC6").Value = -3652
$ws.Range("AD6").Value = -2.3
$ws.Range("AE6").Value = 18865
$ws.Range("AF6").Value = 0.45
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = -3.34
$ws.Range("AJ6").Value = 135786237

# Row 7
$ws.Range("D7").Value = 154310
$ws.Range("E7").Value = 10460
$ws.Range("G7").Value = 4310
$ws.Range("H7").Value = 1825
$ws.Range("I7").Value = -608
$ws.Range("K7").Value = 244058
$ws.Range("L7").Value = 175638
$ws.Range("M7").Value = 68420
$ws.Range("N7").Value = 32025
$ws.Range("P7").Value = 11332
$ws.Range("Q7").Value = -670
$ws.Range("R7").Value = 2118
$ws.Range("S7").Value = -740
$ws.Range("T7").Value = 3100
$ws.Range("U7").Value = -875
$ws.Range("W7").Value = 6.78
$ws.Range("X7").Value = 1.18
$ws.Range("Y7").Value = -2.02
$ws.Range("Z7").Value = 0.74
$ws.Range("AA7").Value = 256.7
$ws.Range("AC7").Value = -330
$ws.Range("AD7").Value = -16.55
$ws.Range("AE7").Value = 16923
$ws.Range("AF7").Value = 0.32
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0

# Row 8
$ws.Range("D8").Value = 159455
$ws.Range("E8").Value = 11305
$ws.Range("G8").Value = 4882
$ws.Range("H8").Value = 3075
$ws.Range("I8").Value = 1035
$ws.Range("K8").Value = 251395
$ws.Range("L8").Value = 180458
$ws.Range("M8").Value = 70938
$ws.Range("N8").Value = 33485
$ws.Range("P8").Value = 11332
$ws.Range("Q8").Value = 9592
$ws.Range("R8").Value = -3728
$ws.Range("S8").Value = -2042
$ws.Range("T8").Value = 2733
$ws.Range("U8").Value = 5550
$ws.Range("W8").Value = 7.09
$ws.Range("X8").Value = 1.93
$ws.Range("Y8").Value = 3.16
$ws.Range("Z8").Value = 1.24
$ws.Range("AA8").Value = 254.39
$ws.Range("AC8").Value = 512
$ws.Range("AD8").Value = 10.68
$ws.Range("AE8").Value = 17695
$ws.Range("AF8").Value = 0.31
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0

# Row 9
$ws.Range("D9").Value = 164542
$ws.Range("E9").Value = 12122
$ws.Range("G9").Value = 5948
$ws.Range("H9").Value = 3852
$ws.Range("I9").Value = 1530
$ws.Range("K9").Value = 260695
$ws.Range("L9").Value = 187042
$ws.Range("M9").Value = 73652
$ws.Range("N9").Value = 34918
$ws.Range("P9").Value = 11332
$ws.Range("Q9").Value = 10148
$ws.Range("R9").Value = -4532
$ws.Range("S9").Value = -1305
$ws.Range("T9").Value = 2913
$ws.Range("U9").Value = 8015
$ws.Range("W9").Value = 7.37
$ws.Range("X9").Value = 2.34
$ws.Range("Y9").Value = 4.47
$ws.Range("Z9").Value = 1.5
$ws.Range("AA9").Value = 253.95
$ws.Range("AC9").Value = 757
$ws.Range("AD9").Value = 7.23
$ws.Range("AE9").Value = 18452
$ws.Range("AF9").Value = 0.3
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0

# Remove cells that no longer have data
$ws.Range("AI7").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("AI9").ClearContents()
